# Applies the "Batterywise analysis" edits to the "Analysis Results" sheet.
# Rows 1-5 and row 11 (Mode) are untouched by the source diff; everything
# from row 6 through the new row 43 gets its label text updated (mostly to
# add measurement units) and, in several spots, its value corrected too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Starting/Ending SoC (%) values were swapped.
$ws.Range("A6").Value = "Starting SoC (%)"
$ws.Range("B6").Value = 99
$ws.Range("A7").Value = "Ending SoC (%)"
$ws.Range("B7").Value = 18

$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"

$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"

# Regenerative Effectiveness: label gains a unit and the sign flips positive.
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("B15").Value = 0.02863166099167976

# Highest/Lowest Cell Voltage rows were swapped (label + value together).
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.443
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.039

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"

# Difference in Temperature: label updated and the previously-blank value is filled in.
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 12

$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# highest/lowest cell temp rows were swapped (label + value together).
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("B28").Value = 47
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("B29").Value = 35

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# From row 31 on, the old "Maximum BMS Temperature in C" row is dropped and
# every following row's label+value shifts up by one, with a brand-new
# "Time spent in 80-90 km/h" row added at the end (row 43).
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 54
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.73834757
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001892140771943574
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 5.573094799614137
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 3.183372796632465
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 3.288608260983952
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 4.937297202490573
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 24.23923528895905
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 6.585986143997194
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 5.371393492940454
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 18.23204419889503
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 28.43549942997457
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
